# Update cryptos list: price and 1h-volume-change values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Cells.Item(2, 4)
$dCell.Value = "'64.945.54"
$dCell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +1.73%  "

$dCell = $ws.Cells.Item(3, 4)
$dCell.Value = "'3.151.77"
$dCell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +3.04%  "

$ws.Cells.Item(4, 5).Value = "  +0.03%  "

$dCell = $ws.Cells.Item(5, 4)
$dCell.Value = "'574.00"
$dCell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.68%  "

$dCell = $ws.Cells.Item(6, 4)
$dCell.Value = "'149.30"
$dCell.Style = "Normal"

$ws.Cells.Item(7, 5).Value = "  +0.04%  "

$dCell = $ws.Cells.Item(8, 4)
$dCell.Value = "'3.151.05"
$dCell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +3.03%  "

$ws.Cells.Item(9, 5).Value = "  +1.84%  "

$dCell = $ws.Cells.Item(10, 4)
$dCell.Value = "'0.159"
$dCell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +3.94%  "

$dCell = $ws.Cells.Item(11, 4)
$dCell.Value = "'6.11"
$dCell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.17%  "

$ws.Cells.Item(12, 5).Value = "  +3.43%  "

$dCell = $ws.Cells.Item(13, 4)
$dCell.Value = "'0.0000264"
$dCell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +13.78%  "

$dCell = $ws.Cells.Item(14, 4)
$dCell.Value = "'37.02"
$dCell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +4.87%  "

$dCell = $ws.Cells.Item(15, 4)
$dCell.Value = "'3.669.25"
$dCell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +3.13%  "

$dCell = $ws.Cells.Item(16, 4)
$dCell.Value = "'65.004.89"
$dCell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +1.79%  "

$dCell = $ws.Cells.Item(17, 4)
$dCell.Value = "'3.150.22"
$dCell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +2.89%  "

$dCell = $ws.Cells.Item(18, 4)
$dCell.Value = "'7.09"
$dCell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +4.57%  "

$ws.Cells.Item(19, 5).Value = "  +1.07%  "

$dCell = $ws.Cells.Item(20, 4)
$dCell.Value = "'505.98"
$dCell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +3.73%  "

$ws.Cells.Item(21, 5).Value = "  +3.07%  "

$ws.Cells.Item(22, 5).Value = "  +3.54%  "

$dCell = $ws.Cells.Item(23, 4)
$dCell.Value = "'15.26"
$dCell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +3.75%  "

$ws.Cells.Item(24, 5).Value = "  +2.75%  "

$dCell = $ws.Cells.Item(25, 4)
$dCell.Value = "'84.17"
$dCell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.87%  "

$ws.Cells.Item(26, 5).Value = "  +0.12%  "

$ws.Cells.Item(27, 5).Value = "  +3.59%  "

$dCell = $ws.Cells.Item(28, 4)
$dCell.Value = "'8.81"
$dCell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +7.63%  "

$ws.Cells.Item(29, 5).Value = "  +5.32%  "

$dCell = $ws.Cells.Item(30, 4)
$dCell.Value = "'2.79"
$dCell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +6.83%  "

$dCell = $ws.Cells.Item(31, 4)
$dCell.Value = "'27.57"
$dCell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +4.07%  "

$dCell = $ws.Cells.Item(32, 4)
$dCell.Value = "'1.00"
$dCell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.03%  "

$ws.Cells.Item(33, 5).Value = "  +2.74%  "

$dCell = $ws.Cells.Item(34, 4)
$dCell.Value = "'6.18"
$dCell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +8.24%  "

$dCell = $ws.Cells.Item(35, 4)
$dCell.Value = "'6.50"
$dCell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +4.05%  "

$dCell = $ws.Cells.Item(36, 4)
$dCell.Value = "'54.87"
$dCell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.33%  "

$dCell = $ws.Cells.Item(37, 4)
$dCell.Value = "'0.0897"
$dCell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +10.19%  "

$dCell = $ws.Cells.Item(38, 4)
$dCell.Value = "'464.44"
$dCell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +5.12%  "

$dCell = $ws.Cells.Item(39, 4)
$dCell.Value = "'0.0419"
$dCell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.90%  "

$dCell = $ws.Cells.Item(40, 4)
$dCell.Value = "'2.98"
$dCell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +8.16%  "

$dCell = $ws.Cells.Item(41, 4)
$dCell.Value = "'8.65"
$dCell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +3.82%  "

$dCell = $ws.Cells.Item(42, 4)
$dCell.Value = "'3.047.14"
$dCell.Style = "Normal"

$ws.Cells.Item(43, 5).Value = "  +0.19%  "

$dCell = $ws.Cells.Item(44, 4)
$dCell.Value = "'2.42"
$dCell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +7.27%  "

$ws.Cells.Item(45, 5).Value = "  +1.64%  "

$dCell = $ws.Cells.Item(46, 4)
$dCell.Value = "'28.47"
$dCell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.35%  "

$ws.Cells.Item(47, 5).Value = "  +12.09%  "

$dCell = $ws.Cells.Item(48, 4)
$dCell.Value = "'0.999"
$dCell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.06%  "

$ws.Cells.Item(49, 5).Value = "  +0.48%  "

$dCell = $ws.Cells.Item(50, 4)
$dCell.Value = "'2.24"
$dCell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +4.24%  "

$dCell = $ws.Cells.Item(51, 4)
$dCell.Value = "'119.62"
$dCell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +1.82%  "
